$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: insert new columns C:G (WIN, TOP4, TOP5, TOP6, RELEGATION)
# and move ExpPoints header from C1 to H1, keeping the same bold/centered style.
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP4"
$ws.Range("E1").Value = "TOP5"
$ws.Range("F1").Value = "TOP6"
$ws.Range("G1").Value = "RELEGATION"
$ws.Range("H1").Value = "ExpPoints"
$ws.Range("B1").Copy()
$ws.Range("C1:H1").PasteSpecial(-4122)

# --- Team name reordering (rows 8-10 and 15-17 swap order) ---
$ws.Range("B8").Value = "Rayo Vallecano"
$ws.Range("B9").Value = "Celta de Vigo"
$ws.Range("B10").Value = "Espanyol"

$ws.Range("B15").Value = "Sevilla"
$ws.Range("B16").Value = "Alavés"
$ws.Range("B17").Value = "Elche"

# --- Updated ExpPoints values, now living in column H ---
$expPoints = @{
    2  = 85.23610550586561
    3  = 83.80602024427213
    4  = 72.03219488286736
    5  = 64.1705632712341
    6  = 60.18625366951277
    7  = 54.83645943015656
    8  = 49.10231464153016
    9  = 48.7699148170532
    10 = 48.30548923886948
    11 = 46.73412709043723
    12 = 46.50935761772428
    13 = 45.21993404140321
    14 = 44.70056806678006
    15 = 44.04239909639792
    16 = 43.53368194546711
    17 = 43.47563047746731
    18 = 41.97133004358815
    19 = 36.6315226125305
    20 = 33.4738955137356
    21 = 33.42406902719281
}

foreach ($row in 2..21) {
    # Clear old ExpPoints value in C and leave new blank prediction columns C:G empty strings
    $ws.Range("C$row").Value = ""
    $ws.Range("D$row").Value = ""
    $ws.Range("E$row").Value = ""
    $ws.Range("F$row").Value = ""
    $ws.Range("G$row").Value = ""
    $ws.Range("H$row").Value = $expPoints[$row]
}
